$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.08750562165990061,0,0.1735998771528813,0.1545093483460214,1.492444062395208,0.002468405472996795,0,0.5915225028457844,0.1779213817055343,0.6615644409347112,0,0.2844807193773704,0,3.648149885222466),
    @(0.07660246508839919,0,0.1710091665664635,0.1540234438458796,1.497720545593317,0.002471128877247946,0,0.6005887154299057,0.1785819278493221,0.5857526274613178,0,0.2662373191513439,0,3.669841123677429),
    @(0.06989396031504214,0,0.1694796270467833,0.1537931766480831,1.501789723700625,0.002472890959814052,0,0.6064735899519835,0.1790758196068722,0.5391368886814973,0,0.2551044755012057,0,3.685296547345203),
    @(0.06715689037180539,0,0.1688717832744118,0.1537164904591606,1.503656530966374,0.002473631696406371,0,0.608951727727773,0.1792993007679016,0.5201248476128058,0,0.2505852932414427,0,3.692131925044833),
    @(0.0667022086343394,0,0.1687717869932897,0.153704793458779,1.503979111634592,0.00247375606649558,0,0.60936805246482,0.1793377517196362,0.5169669937786523,0,0.2498359539994297,0,3.69329937270237),
    @(0.0698570603257167,0,0.1694713667719867,0.1537920729529638,1.501814055534595,0.002472900857655484,0,0.6065066870608917,0.1790787435870307,0.5388805479607299,0,0.2550434568403261,0,3.685386556912263),
    @(0.08374921996278317,0,0.1726939487032695,0.1543276960869129,1.494091265418888,0.002469325884413756,0,0.5945824912421003,0.1781308136606654,0.6354390603684408,0,0.278176297993582,0,3.655185396293575),
    @(0.1108737213141637,0,0.1794959217301653,0.1559171271519055,1.485528025431734,0.002463025624448824,0,0.5737237704711289,0.1769723915937931,0.8242198151533273,0,0.3240751160357078,0,3.612926910411801),
    @(0.1307220213705591,0,0.1847843432112768,0.157412350861005,1.483251013175192,0.002458825503840918,0,0.5599373846685278,0.17654816402937,0.9625310382385237,0,0.3581141133362138,0,3.592240408156187),
    @(0.1397325443067388,0,0.1872527532507462,0.1581634478461886,1.483087484751707,0.002457006927427868,0,0.553999549908756,0.1764478434616237,1.025361342088274,0,0.3736664675152923,0,3.585083462719979),
    @(0.14314174896829,0,0.1881964231758246,0.1584580405690232,1.483151016288417,0.002456331451461432,0,0.5517990342603896,0.1764231748524736,1.049139956254862,0,0.3795652711904651,0,3.582697645872486),
    @(0.1424076474561673,0,0.187992790736061,0.1583941429582083,1.483131753605875,0.002456476342070912,0,0.5522708187279868,0.1764278953072917,1.044019441323144,0,0.3782944419571166,0,3.583197042998989),
    @(0.1400130806685809,0,0.1873302108796935,0.1581874805361956,1.483090197081253,0.0024569510918849,0,0.5538175497070501,0.176445547018929,1.027317907467591,0,0.3741515776216957,0,3.584880677584124),
    @(0.138545957326329,0,0.1869255231206637,0.1580622172410848,1.483081081228562,0.002457243604257208,0,0.5547712208669209,0.176458093814972,1.017085896043454,0,0.3716151775500691,0,3.585954203505196),
    @(0.1301327704832289,0,0.1846242815370118,0.1573646894061156,1.483279254781834,0.002458946199792991,0,0.5603321531009335,0.1765565844686492,0.9584230549853316,0,0.3570990706902819,0,3.59275350388404),
    @(0.1249666439389046,0,0.1832285436870933,0.1569549207197163,1.483624258925929,0.002460014226705426,0,0.563829087227258,0.1766407354695119,0.9224118792510296,0,0.3482110782036756,0,3.597502044569922),
    @(0.1219934896173527,0,0.1824316565665924,0.1567259089077382,1.483904799563675,0.002460637198391994,0,0.5658718350033451,0.1766978590005408,0.901690975080669,0,0.3431053431035807,0,3.60044536570183),
    @(0.1209865392268341,0,0.1821628607026895,0.156649516972017,1.484013886489564,0.002460849616918455,0,0.5665688669153064,0.1767186982021514,0.8946738557930587,0,0.3413777375534792,0,3.601478339905071),
    @(0.1255167680511988,0,0.1833765118403363,0.1569978505306544,1.483579035097634,0.002459899636422917,0,0.5634535822075963,0.176630874777679,0.9262461905550765,0,0.3491565595925934,0,3.596974601861604),
    @(0.1407165029784352,0,0.187524584718048,0.1582479065535374,1.483098998228755,0.002456811289019769,0,0.553361933984251,0.1764400007996443,1.032223938300319,0,0.3753681827483391,0,3.58437734762947),
    @(0.1506334946015926,0,0.1902876440453127,0.1591241475567138,1.483516530587792,0.002454869667994658,0,0.5470463337311076,0.1763928943998323,1.101405202334661,0,0.3925540266959544,0,3.578035097709886),
    @(0.1453422278120797,0,0.1888082096732404,0.1586510686274352,1.483226768140852,0.002455898941378784,0,0.5503914687031619,0.1764109332163599,1.064489689043398,0,0.3833766787639803,0,3.581246963548011),
    @(0.125268066516071,0,0.1833095981730253,0.1569784214975627,1.483599224763722,0.002459951415046644,0,0.5636232473743448,0.1766353055593868,0.9245127531475532,0,0.3487290946948463,0,3.597212394241893),
    @(0.1035493023219942,0,0.1776044774058789,0.1554295527795162,1.487139807927058,0.002464654421666154,0,0.5790963900795574,0.1772108020837209,0.7732146541736995,0,0.311601913206097,0,3.62254077838648)
)

$startRow = 2
$startCol = 2  # column B
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowVals[$c]
    }
}

Write-Host "Updated B2:O25"